$d = $word.ActiveDocument
$d.Content.Find.Execute("19.01.2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "22.01.2024", 2)
